# Update "想去人数" (F column) figures for the Suzhou con-info workbook.
# The same set of events (by row) appears on both the "展览" (Exhibition)
# sheet and the "全部类型" (All types) sheet; the "全部类型" sheet has one
# extra row inserted above row 33, so the row numbers differ by one from
# that point on.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value, for the "展览" worksheet.
$exhibitionUpdates = @{
    4  = 1139
    6  = 89
    8  = 66
    9  = 1174
    10 = 16569
    11 = 288
    12 = 207
    13 = 1044
    14 = 6422
    18 = 32
    19 = 129
    21 = 57
    24 = 40
    26 = 10
    28 = 231
    29 = 905
    30 = 67
    31 = 5069
    33 = 11416
    35 = 23
    36 = 159
    37 = 216
    40 = 76
}

# Row -> new F-column value, for the "全部类型" worksheet.
$allTypesUpdates = @{
    4  = 1139
    6  = 89
    8  = 66
    9  = 1174
    10 = 16569
    11 = 288
    12 = 207
    13 = 1044
    14 = 6422
    18 = 32
    19 = 129
    21 = 57
    24 = 40
    26 = 10
    28 = 231
    29 = 905
    30 = 67
    31 = 5069
    34 = 11416
    36 = 23
    37 = 159
    38 = 216
    41 = 76
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
